$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("E3").Value = "['MCT-3A-Eletropneumática', -, -, -]"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "['MCT-3A-Eletropneumática', -, -, -]"
$ws.Range("E6").Value = "['MCT-3A-Eletropneumática', -, -, -]"
$ws.Range("E7").Value = "['MCT-3A-Eletropneumática', -, -, -]"
$ws.Range("C8").Value = "-"

$ws.Range("B11").Value = "['MEC-3A-C.pneumática', -, -, -]"
$ws.Range("D11").Value = "-"
$ws.Range("B12").Value = "['MEC-3A-C.pneumática', -, -, -]"
$ws.Range("D12").Value = "-"
$ws.Range("B14").Value = "['MEC-3A-C.pneumática', -, -, -]"
$ws.Range("D14").Value = "-"
$ws.Range("B15").Value = "['MEC-3A-C.pneumática', -, -, -]"
$ws.Range("D15").Value = "-"

$ws.Range("D19").Value = "[-, 'MEC-1NB-Desenho tecnico mecanico']"
$ws.Range("B20").Value = "['MEC-1NA-Desenho tecnico mecanico – T1', 'MEC-1NA-Desenho tecnico mecanico – T1']"
$ws.Range("E20").Value = "-"
$ws.Range("B21").Value = "['MEC-1NA-Desenho tecnico mecanico – T1', 'MEC-1NA-Desenho tecnico mecanico – T1']"
$ws.Range("E21").Value = "-"
